# Append the new run-log row (row 20) to Sheet1, mirroring the format of
# the existing rows (row 19 is the current last row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 20
$lastRow = $newRow - 1

# Copy row 19's formatting (styles) down into the new row first so the
# appended row matches the look of every other data row.
$ws.Range("A$lastRow`:H$lastRow").Copy()
$ws.Range("A$newRow`:H$newRow").PasteSpecial(-4122)  # xlPasteFormats

# Fill in the new row's values.
$ws.Cells.Item($newRow, 1).Value = "2025-08-16 09:37:41 UTC"
$ws.Cells.Item($newRow, 2).Value = "2025-08-16 15:07:41 IST"
$ws.Cells.Item($newRow, 3).Value = "SKIPPED"
$ws.Cells.Item($newRow, 4).Value = "No change in PDF. Skipping download & Excel update."
$ws.Cells.Item($newRow, 5).Value = "https://nalcoindia.com/wp-content/uploads/2025/08/INGOT-15-08-2025.pdf"
$ws.Cells.Item($newRow, 7).Value = 0
# Columns F (Saved PDF) and H (Total Rows After) stay blank for a SKIPPED run.

$wb.Save()
